$d = $word.ActiveDocument

# 1. Insert a new bullet paragraph before "Each time when player starts to play..."
#    containing the teleport-to-tutorial-room requirement.
$insRng = $d.Content
$insRng.Find.Execute("Each time when player starts to play, a map should be generated following BSP.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$insRng.Collapse(1)
$insRng.InsertParagraphBefore()
$insRng.Text = "When pressing T, the player should be teleported to the tutorial room, where one can try out key mappings we have."

# 2. Delete the run of bullet paragraphs from "When starting a new game..." through
#    "...sound track...different from the main menu." (7 paragraphs).
$startRng = $d.Content
$startRng.Find.Execute("When starting a new game, or seeing a boss, the player should be able to interact with the dialogue box.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$startPos = $startRng.Paragraphs(1).Range.Start

$endRng = $d.Content
$endRng.Find.Execute("When players enter the main menu, they should hear a particular sound track.", `
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$endPos = $endRng.Paragraphs(1).Range.End

$delRng = $d.Range($startPos, $endPos)
$delRng.Delete()
